# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets,
# reflecting refreshed scrape counts (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14
$ws1.Range("F3").Value = 1077
$ws1.Range("F4").Value = 182
$ws1.Range("F5").Value = 3075
$ws1.Range("F6").Value = 98
$ws1.Range("F7").Value = 294
$ws1.Range("F9").Value = 8
$ws1.Range("F10").Value = 8
$ws1.Range("F11").Value = 129
$ws1.Range("F12").Value = 107
$ws1.Range("F13").Value = 171
$ws1.Range("F14").Value = 79
$ws1.Range("F15").Value = 2784
$ws1.Range("F16").Value = 1070
$ws1.Range("F17").Value = 4

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14
$ws4.Range("F4").Value = 1077
$ws4.Range("F5").Value = 182
$ws4.Range("F6").Value = 3075
$ws4.Range("F7").Value = 98
$ws4.Range("F8").Value = 294
$ws4.Range("F11").Value = 8
$ws4.Range("F12").Value = 8
$ws4.Range("F13").Value = 129
$ws4.Range("F14").Value = 107
$ws4.Range("F15").Value = 171
$ws4.Range("F16").Value = 79
$ws4.Range("F17").Value = 2784
$ws4.Range("F18").Value = 1070
$ws4.Range("F19").Value = 4

$wb.Save()
